$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (header is row 1, data starts row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Add new header cells (Wins / Losses / Ties), reusing the formatting
# of the neighboring header cell (AC1) so the new cells pick up the
# existing bold/border/centered header style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 98   # AD = Wins
    $ws.Cells.Item($r, 31).Value = 64   # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = Ties
}
